# re-run RU 1001; without crop
# Update mean values in the country-comparison sheet (rows 2-9) to reflect
# a re-run of the RU (Russia, column L) computation without cropping.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 67.8255122017956
$ws.Range("L2").Value = 73.9047667329172

$ws.Range("B3").Value = 55.3746767090015
$ws.Range("D3").Value = 62.2273519140294
$ws.Range("E3").Value = 56.7933697165467
$ws.Range("F3").Value = 75.7304011850663
$ws.Range("G3").Value = 53.058303552895
$ws.Range("H3").Value = 68.9059343539561
$ws.Range("I3").Value = 59.3195783835747
$ws.Range("J3").Value = 62.3575700142156
$ws.Range("K3").Value = 54.8090076208069
$ws.Range("L3").Value = 49.0485989036895
$ws.Range("M3").Value = 84.6428606188793
$ws.Range("N3").Value = 47.4912333237318

$ws.Range("B4").Value = 40.563158268296
$ws.Range("L4").Value = 36.9618773883661

$ws.Range("B6").Value = 64.7210522905015

$ws.Range("B7").Value = 67.1354010141054
$ws.Range("C7").Value = 73.1501666372061
$ws.Range("L7").Value = 63.442019211072

$ws.Range("B8").Value = 68.4899159160604
$ws.Range("L8").Value = 60.1821932205212

$ws.Range("B9").Value = 61.8798692282585
$ws.Range("L9").Value = 54.0183622108344
